# Added division and address supplement
# Insert a new "Abteilung" (division) column right before the existing
# "Anrede" column, and a new "Adresszusatz" (address supplement) column
# right before the existing "Postleitzahl_2" column. Excel shifts the
# existing columns (and their explicit widths) to the right automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert "Abteilung" before "Anrede" (was column N, now becomes O)
$ws.Columns("N").Insert()
$ws.Range("N1").Value = "Abteilung"

# Insert "Adresszusatz" before "Postleitzahl_2" (was column S, now becomes U)
$ws.Columns("T").Insert()
$ws.Range("T1").Value = "Adresszusatz"

# Give the new "Adresszusatz" column an explicit width (matches the
# neighbouring "Hausnummer" column's width)
$ws.Columns("T").ColumnWidth = 12

# Leave the cursor where the user ended up after typing the new header
$null = $ws.Range("T2").Select()
